$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27 ("VERIFY_TEXT_PRESENT ..."), shifting the
# remaining rows (old 27-29) down to 28-30. The new row becomes a "WAIT"
# step, matching the other WAIT rows already used between CLICK_JS steps
# (e.g. rows 20/22/24/26).
$ws.Rows.Item(27).Insert()

$newRow = $ws.Range("A27:E27")
$newRow.Borders.LineStyle = 1

$ws.Cells.Item(27, 2).Value = "WAIT"

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("B26").Select()
